$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    "52-27=25",
    "71-52=19",
    "9+22=31",
    "81-73=8",
    "25+66=91",
    "7+78=85",
    "90-31=59",
    "84-68=16",
    "41-5=36",
    "71-22=49",
    "19+33=52",
    "90-62=28",
    "19+15=34",
    "91-52=39",
    "18+58=76",
    "48-9=39",
    "59+19=78",
    "61-5=56",
    "78+7=85",
    "91-8=83",
    "54-25=29",
    "14+57=71",
    "6+19=25",
    "19+18=37",
    "56-37=19",
    "85+8=93",
    "74-37=37",
    "49+13=62",
    "68+26=94",
    "90-18=72",
    "71-68=3",
    "6+29=35",
    "73-34=39",
    "59+23=82",
    "81-14=67",
    "7+25=32",
    "35-26=9",
    "55+38=93",
    "75-29=46",
    "91-84=7",
    "19+6=25",
    "63-35=28",
    "34-5=29",
    "38+7=45",
    "35+6=41",
    "28+69=97",
    "11-8=3",
    "9+37=46",
    "19+57=76",
    "76-38=38",
    "18+69=87",
    "25+17=42",
    "75+16=91",
    "17+16=33",
    "60-11=49",
    "15+77=92",
    "61-48=13",
    "5+37=42",
    "88-39=49",
    "19+27=46",
    "27+8=35",
    "59+23=82",
    "96-57=39",
    "92-18=74",
    "54+29=83",
    "95-79=16",
    "59+7=66",
    "47+16=63",
    "96-38=58",
    "52-46=6",
    "35+9=44",
    "94-48=46",
    "64+29=93",
    "65+29=94",
    "71-64=7",
    "46+7=53",
    "64-45=19",
    "6+27=33",
    "60-5=55",
    "75-38=37",
    "91-58=33",
    "36+18=54",
    "24+49=73",
    "26+59=85",
    "93-48=45",
    "17+78=95",
    "10-6=4",
    "64-49=15",
    "49+26=75",
    "82-19=63",
    "28+58=86",
    "27+27=54",
    "91-87=4",
    "94-17=77",
    "53-28=25",
    "87+7=94",
    "6+55=61",
    "51-23=28",
    "92-88=4",
    "50-5=45"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

if ($rows * $cols -ne $values.Count) {
    throw "Table shape $rows x $cols ($($rows * $cols) cells) does not match expected $($values.Count) values"
}

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$i]
        $i++
    }
}

Write-Output ("Updated " + $i + " cells")
